# Append a new bullet to the change-log list, right after the last
# entry ("Erstellung der Buttons reorganisiert"). InsertParagraphAfter
# clones the paragraph's style/numbering (Listenabsatz, ilvl 0, numId 1)
# onto the freshly inserted paragraph.
$d = $word.ActiveDocument

$lastParagraph = $d.Paragraphs.Last
$lastParagraph.Range.InsertParagraphAfter() | Out-Null

$newParagraph = $d.Paragraphs.Last
$newParagraph.Range.Text = "Laden und Speichern von Bildern"
